# Generate Report for Archive
#
# The two "e2e" entries af22f8cb-2408-4abf-8703-dc848a9921b7 and
# dbe99156-d094-4e36-8ed0-b478e193b6d9 swap report rows (row 5 <-> row 6)
# across all three sheets, and dbe99156's Status flips from
# "Ready for handoff" to "In Translation".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview": columns A (File Name), B (Path And Name, hyperlinked),
# E (zh-cn), F (de-de), G (Latest HO Xliff Generate Date)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A5").Value = "dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
$wsOverview.Range("B5").Value = "e2e\dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
$wsOverview.Range("E5").Value = "In Translation"
$wsOverview.Range("F5").Value = "In Translation"
$wsOverview.Range("G5").Value = "2016-09-02 12:47:53"

$wsOverview.Range("A6").Value = "af22f8cb-2408-4abf-8703-dc848a9921b7.md"
$wsOverview.Range("B6").Value = "e2e\af22f8cb-2408-4abf-8703-dc848a9921b7.md"
$wsOverview.Range("E6").Value = "Ready for handoff"
$wsOverview.Range("F6").Value = "Ready for handoff"
$wsOverview.Range("G6").Value = "2016-09-02 12:48:22"

foreach ($h in $wsOverview.Hyperlinks) {
  $addr = $h.Range.Address()
  if ($addr -eq '$B$5') {
    $h.TextToDisplay = "e2e\dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
  }
  if ($addr -eq '$B$6') {
    $h.TextToDisplay = "e2e\af22f8cb-2408-4abf-8703-dc848a9921b7.md"
  }
}

# ---------------------------------------------------------------------
# Sheet "zh-cn": columns A (Source File Name), C (Status),
# G (Latest Handoff File), H (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A5").Value = "dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
$wsZhCn.Range("C5").Value = "In Translation"
$wsZhCn.Range("G5").Value = "dbe99156-d094-4e36-8ed0-b478e193b6d9.e834e5271e318c2dd5dd1d279dc98e8633dbf018.zh-cn.xlf"
$wsZhCn.Range("H5").Value = "2016-09-02 12:47:48"

$wsZhCn.Range("A6").Value = "af22f8cb-2408-4abf-8703-dc848a9921b7.md"
$wsZhCn.Range("C6").Value = "Ready for handoff"
$wsZhCn.Range("G6").Value = "af22f8cb-2408-4abf-8703-dc848a9921b7.8f1f19c580ec54f667705e4257c60c3de87db8ea.zh-cn.xlf"
$wsZhCn.Range("H6").Value = "2016-09-02 12:48:18"

foreach ($h in $wsZhCn.Hyperlinks) {
  $addr = $h.Range.Address()
  if ($addr -eq '$A$5') {
    $h.TextToDisplay = "dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
  }
  if ($addr -eq '$A$6') {
    $h.TextToDisplay = "af22f8cb-2408-4abf-8703-dc848a9921b7.md"
  }
}

# ---------------------------------------------------------------------
# Sheet "de-de": columns A (Source File Name), C (Status),
# G (Latest Handoff File), H (Latest Handoff Datetime)
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A5").Value = "dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
$wsDeDe.Range("C5").Value = "In Translation"
$wsDeDe.Range("G5").Value = "dbe99156-d094-4e36-8ed0-b478e193b6d9.e834e5271e318c2dd5dd1d279dc98e8633dbf018.de-de.xlf"
$wsDeDe.Range("H5").Value = "2016-09-02 12:47:53"

$wsDeDe.Range("A6").Value = "af22f8cb-2408-4abf-8703-dc848a9921b7.md"
$wsDeDe.Range("C6").Value = "Ready for handoff"
$wsDeDe.Range("G6").Value = "af22f8cb-2408-4abf-8703-dc848a9921b7.8f1f19c580ec54f667705e4257c60c3de87db8ea.de-de.xlf"
$wsDeDe.Range("H6").Value = "2016-09-02 12:48:22"

foreach ($h in $wsDeDe.Hyperlinks) {
  $addr = $h.Range.Address()
  if ($addr -eq '$A$5') {
    $h.TextToDisplay = "dbe99156-d094-4e36-8ed0-b478e193b6d9.md"
  }
  if ($addr -eq '$A$6') {
    $h.TextToDisplay = "af22f8cb-2408-4abf-8703-dc848a9921b7.md"
  }
}
